$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells
$ws.Range("A1").Value = "NOMOR"
$ws.Range("B1").Value = "ITEM_ASSY"

$rng = $ws.Range("B1")
$rng.Font.Name = "Calibri"
$rng.Font.Bold = $true
$rng.Interior.Color = 52479
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

Write-Host "done"
